# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share that string).
# 2) Narrow the "Status" columns (Overview E & F, zh-cn C, de-de C) from
#    their old ~17.22-char width down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status value wherever it shows up ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
